{"js": "// Commit: \"Desi demonstraing how git works\"\n// The document originally contains:\n//   1. A paragraph with the text \"Actions\"\n//   2-5. A 3-level bulleted list (\"Possible Suspension\" > \"Components\" > \"Springs\"/\"Drive belts\")\n//   6. An otherwise-empty paragraph (same list formatting) that only carries the\n//      \"_GoBack\" bookmark.\n// After the edit the whole list is gone and the document is reduced to a single,\n// un-styled paragraph whose text reads \"Showing Nadav how git works\" and which still\n// carries the \"_GoBack\" bookmark at its end.\n\nconst body = context.document.body;\n\n// Load all paragraphs currently in the body.\nlet paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Remove every paragraph except the first one (this drops the bulleted list items\n// together with the trailing empty/bookmarked paragraph - and, as a side effect,\n// the \"_GoBack\" bookmark that lived on that last paragraph).\nfor (let i = paragraphs.items.length - 1; i >= 1; i--) {\n  paragraphs.items[i].delete();\n}\nawait context.sync();\n\n// Re-fetch the (now only) paragraph and swap its text.\nparagraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst firstParagraph = paragraphs.items[0];\nconst contentRange = firstParagraph.getRange(\"Content\");\ncontentRange.insertText(\"Showing Nadav how git works\", Word.InsertLocation.replace);\nawait context.sync();\n\n// Restore the \"_GoBack\" bookmark at the end of the (now only) paragraph, exactly how\n// it appeared in the original document.\nconst endRange = firstParagraph.getRange(\"End\");\nendRange.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Commit: \"Desi demonstraing how git works\"\n#\n# The document originally contains:\n#   1. A paragraph with the text \"Actions\"\n#   2-5. A 3-level bulleted list (\"Possible Suspension\" > \"Components\" > \"Springs\"/\"Drive belts\")\n#   6. An otherwise-empty paragraph (same list formatting) that only carries the\n#      \"_GoBack\" bookmark.\n#\n# After the edit the whole list is gone and the document is reduced to a single,\n# un-styled paragraph whose text reads \"Showing Nadav how git works\" and which still\n# carries the \"_GoBack\" bookmark at its end.\n\n$d = $word.ActiveDocument\n\n# Remove every paragraph after the first one. This drops the bulleted list items\n# together with the trailing, otherwise-empty paragraph that only carried the\n# \"_GoBack\" bookmark (that bookmark is re-created further down).\nfor ($i = $d.Paragraphs.Count; $i -ge 2; $i--) {\n    $d.Paragraphs.Item($i).Range.Delete()\n}\n\n# Swap the text of the (now only) paragraph, keeping its paragraph mark untouched.\n$p1 = $d.Paragraphs.Item(1)\n$bodyRange = $p1.Range.Duplicate\n$bodyRange.Collapse(1)   # wdCollapseStart\n$bodyRange.Text = \"Showing Nadav how git works\"\n\n# Re-create the \"_GoBack\" bookmark collapsed at the end of the paragraph's text (i.e.\n# right before the paragraph mark), exactly as it originally appeared. A temporary\n# placeholder character is inserted first and removed afterwards so the collapsed\n# bookmark lands right after the run instead of being normalized into one that spans\n# (and so appears to wrap) the whole paragraph.\n$p1 = $d.Paragraphs.Item(1)\n$endPos = $p1.Range.End - 1\n\n$placeholder = $d.Range($endPos, $endPos)\n$placeholder.InsertAfter(\"X\")\n\n$bookmarkRange = $d.Range($endPos, $endPos)\n$d.Bookmarks.Add(\"_GoBack\", $bookmarkRange)\n\n$d.Range($endPos, $endPos + 1).Delete()\n"}
